$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EmployeeData")

# Update FirstName column (A) test data values
$ws.Range("A2").Value = "Alida74"
$ws.Range("A3").Value = "S954"
$ws.Range("A4").Value = "S1122"
$ws.Range("A5").Value = "Sa522"

# Update Username column (E) test data values
$ws.Range("E2").Value = "alippli236823"
$ws.Range("E3").Value = "desırrı123523"
$ws.Range("E4").Value = "sde234223"
$ws.Range("E5").Value = "Saa1231123"
